$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "phi" block, right after "gamma_q"/"gamma_w".
# This shifts the old O..V columns one slot to the right (P..W) and leaves a
# blank column O that we populate with the new "gamma_wc" data.
$ws.Columns("O").EntireColumn.Insert()

# "gamma_w" (still sitting in N1) becomes "gamma_wf" - its value (1.4) is unchanged.
$ws.Range("N1").Value = "gamma_wf"

# The freshly inserted column gets the new "gamma_wc" header and value.
$ws.Range("O1").Value = "gamma_wc"
$ws.Range("O2").Value = 1.8

# The new column wasn't produced by AutoFit/BestFit, so give it an explicit
# (non bestFit) width of 13 characters, matching the other data columns.
$ws.Columns("O").ColumnWidth = 12.166666666666666

# Match the author's final cursor position/selection on the sheet.
[void]$ws.Range("O3").Select()

Write-Host "done"
